# NMCARS-PART-5208.docx fix:
#   - The paragraph "(5)(iii) The HCA is the designated management
#     official. ..." should become two runs within the SAME paragraph:
#     one run holding "(5)" and a second run holding the remainder
#     "(iii) The HCA is the designated management official. ...".
#   - The paragraph's style changes from "List1" to "List2".

$d = $word.ActiveDocument

# Locate the target paragraph via Find; Find.Execute collapses the
# range to the matched text, so Paragraphs(1) is the paragraph we want.
$target = $d.Content
$found = $target.Find.Execute(
    "(5)(iii) The HCA is the designated management official. The management official shall coordinate any request for waiver with the DON CIO.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $target.Paragraphs(1)

    # Change the paragraph style List1 -> List2.
    $para.Style = "List2"

    # Split the single run into two runs ("(5)" | "(iii) ... DON CIO.")
    # without altering any visible formatting. Inserting then removing a
    # bookmark at the split point forces the text either side of it to be
    # stored as separate runs, while leaving no visible trace behind.
    $paraStart = $para.Range.Start
    $splitPoint = $d.Range($paraStart + 3, $paraStart + 3)
    $d.Bookmarks.Add("TempRunSplitMark", $splitPoint)
    $d.Bookmarks("TempRunSplitMark").Delete()
}
